$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sending/target cluster vocabulary changed: a new "ECs" cluster was
# added, so the 2 original rows (FAPs->MuSCs, MuSCs->MuSCs) are replaced by
# 4 rows covering every combination of the two sending clusters (ECs, FAPs)
# against the two target clusters (ECs, MuSCs), each with refreshed TPM
# based values.
# NOTE: PowerShell hashtable keys are case-insensitive, so the row-number
# key is named "Row" (not "r") to avoid colliding with the "R" column.

$rows = @(
  @{ Row = 2; A = "ECs";  D = "ECs";   E = 1; F = 0.3333333333333333; G = 0.06729733333333333; H = 0.201892;   I = 0.01373511018321553; J = 0.01373511018321553; K = 1; L = 0.3333333333333333; M = 0.044174; N = 0.132522; O = 0.2474193313505733; P = 0.2474193313505733; Q = 0.002972792402666667;  R = 0.026755131624;    S = 0.003398331777557637; T = 0.003398331777557636 },
  @{ Row = 3; A = "ECs";  D = "MuSCs"; E = 1; F = 0.3333333333333333; G = 0.06729733333333333; H = 0.201892;   I = 0.01373511018321553; J = 0.01373511018321553; K = 2; L = 0.6666666666666666; M = 0.134365; N = 0.403095; O = 0.7525806686494267; P = 0.7525806686494266; Q = 0.009042406193333333;  R = 0.08138165573999999; S = 0.0103367784056579;   T = 0.01033677840565789 },
  @{ Row = 4; A = "FAPs"; D = "ECs";   E = 3; F = 1;                   G = 4.83236;             H = 14.49708;  I = 0.9862648898167845;  J = 0.9862648898167844;  K = 1; L = 0.3333333333333333; M = 0.044174; N = 0.132522; O = 0.2474193313505733; P = 0.2474193313505733; Q = 0.21346467064;         R = 1.92118203576;      S = 0.2440209995730157;   T = 0.2440209995730156 },
  @{ Row = 5; A = "FAPs"; D = "MuSCs"; E = 3; F = 1;                   G = 4.83236;             H = 14.49708;  I = 0.9862648898167845;  J = 0.9862648898167844;  K = 2; L = 0.6666666666666666; M = 0.134365; N = 0.403095; O = 0.7525806686494267; P = 0.7525806686494266; Q = 0.6493000514;          R = 5.8437004626;       S = 0.7422438902437689;   T = 0.7422438902437687 }
)

foreach ($row in $rows) {
  $r = $row.Row
  $ws.Range("A$r").Value = $row.A
  $ws.Range("B$r").Value = "Rspo3"
  $ws.Range("C$r").Value = "Lgr6"
  $ws.Range("D$r").Value = $row.D
  $ws.Range("E$r").Value = $row.E
  $ws.Range("F$r").Value = $row.F
  $ws.Range("G$r").Value = $row.G
  $ws.Range("H$r").Value = $row.H
  $ws.Range("I$r").Value = $row.I
  $ws.Range("J$r").Value = $row.J
  $ws.Range("K$r").Value = $row.K
  $ws.Range("L$r").Value = $row.L
  $ws.Range("M$r").Value = $row.M
  $ws.Range("N$r").Value = $row.N
  $ws.Range("O$r").Value = $row.O
  $ws.Range("P$r").Value = $row.P
  $ws.Range("Q$r").Value = $row.Q
  $ws.Range("R$r").Value = $row.R
  $ws.Range("S$r").Value = $row.S
  $ws.Range("T$r").Value = $row.T
}
